# edit.ps1 - apply the "new solutions of tasks" wording fixes
# (01. Print an Array with a Given Delimiter - exercise statement)
#
# We locate each target phrase with Find (no replacement performed by
# Find itself), then overwrite the located Range's .Text directly. Doing
# the substitution through Find.Execute's own ReplaceWith argument runs
# the text through the "smart quotes" autocorrect pass and mangles plain
# apostrophes ('), so we avoid that path entirely.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
    if (-not $found) {
        Write-Host "WARNING: not found -> $old"
        return $false
    }
    $range.Text = $new
    return $true
}

# 1. "resulted" -> "resulting"
Replace-Text " is the resulted array after the rotations" " is the resulting array after the rotations" | Out-Null

# 2. Heading: "Subsequence" -> "Subset"
Replace-Text "Extract Increasing Subsequence from Array" "Extract Increasing Subset from Array" | Out-Null

# 3. Bold run "non-decreasing subsequence" -> "non-decreasing subset"
#    (this is the first of two occurrences in the doc; Replace-Text with
#     wdReplaceOne only touches the first match, which is this one)
Replace-Text "non-decreasing subsequence" "non-decreasing subset" | Out-Null

# 4. "which should be a non-decreasing subsequence." -> "...subset."
Replace-Text "which should be a non-decreasing subsequence." "which should be a non-decreasing subset." | Out-Null

# 5. " one and so on. " -> " one, and so on. "
Replace-Text " one and so on. " " one, and so on. " | Out-Null

# 6. "(the result of the compare is 0)" -> "(the result of the comparison is 0)"
Replace-Text " (the result of the compare is 0), we need to compare by the " " (the result of the comparison is 0), we need to compare by the " | Out-Null

# 7. "will be first player's" -> "will be the first player's" and add comma before "and so on."
Replace-Text "so the first element of the input will be first player's chosen coordinates, the second element will be the second player's turn coordinates and so on." "so the first element of the input will be the first player's chosen coordinates, the second element will be the second player's turn coordinates, and so on." | Out-Null

# 8. "in row that players make" -> "in a row that players make"
Replace-Text "in row that players make" "in a row that players make" | Out-Null

# 9. "each row the dashboard" -> "each row of the dashboard"
Replace-Text "the elements of each row the dashboard should be separated by" "the elements of each row of the dashboard should be separated by" | Out-Null

# 10. "on new line." -> "on a new line."
Replace-Text "and each row should be on new line." "and each row should be on a new line." | Out-Null

# 11. "...diagonal sum or the original matrix, if the two diagonals..." -> "...diagonal sum, or the original matrix if the two diagonals..."
Replace-Text "l are changed to the diagonal sum or the original matrix, if the two diagonals have different sums." "l are changed to the diagonal sum, or the original matrix if the two diagonals have different sums." | Out-Null

# 12. Remove the stray "_GoBack" bookmark after "Examples"
try {
    $gb = $d.Bookmarks("_GoBack")
    $gb.Delete()
} catch {
    Write-Host "no _GoBack bookmark"
}
